$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) interpretation for price cells that look like plain numbers,
# so trailing zeros / exact formatting are preserved exactly as text (matching source data).
$textRefs = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

# Apply the updated cell values from the latest cryptos snapshot.
$ws.Range("D2").Value = "29.456.70"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.851.12"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "240.21"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "0.6290"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07658"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "0.2916"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "24.82"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").Value = "2.163.71"
$ws.Range("E11").Value = "  +16.78%  "
$ws.Range("D12").Value = "0.07744"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "5.033"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "0.6812"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "0.00001060"
$ws.Range("E15").Value = "  -4.94%  "
$ws.Range("D16").Value = "83.53"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "29.573.67"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "229.11"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "7.458"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "157.30"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "0.1388"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "8.431"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "17.74"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").Value = "1.395"
$ws.Range("E28").Value = "  +7.34%  "
$ws.Range("D29").Value = "1.462"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "0.05614"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "4.133"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "4.065"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "1.846"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "0.7004"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").Value = "2.595"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.232.54"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01805"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "6.465"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").Value = "0.9084"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "102.65"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "66.08"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "7.211"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "0.00000000118"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").Value = "0.4029"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "0.1154"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "1.681"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +0.11%  "
